# Split the single run in the "hamburger navigation" paragraph into several
# runs, changing "So, to fix this issue a created" to "To fix this issue, I created"
# (exactly matching the author's commit) while leaving all other text, and all
# paragraph/run formatting, untouched.

$d = $word.ActiveDocument

# Locate the old phrase (including the leading space that separates it from
# "...understand what that is for.") inside the document text.
$oldPhrase = " So, to fix this issue a created"
$content = $d.Content.Text
$startPos = $content.IndexOf($oldPhrase)

if ($startPos -lt 0) {
    throw "Could not find target phrase to edit."
}

# Find the paragraph that contains the match, so we know exactly where its
# text (excluding the trailing paragraph mark) ends.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End
    if ($startPos -ge $pStart -and $startPos -lt $pEnd) {
        $targetPara = $para
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find paragraph containing target phrase."
}

$paraTextEnd = $targetPara.Range.End - 1

# Replace everything from the start of the old phrase through the end of the
# paragraph's text with freshly split runs. (The tail of the paragraph after
# the edited phrase is reproduced unchanged as the final run.)
$editRange = $d.Range($startPos, $paraTextEnd)

$tailAfterPhrase = $editRange.Text.Substring($oldPhrase.Length)

# Escape the bits of text we splice into the XML fragment so this keeps working
# even if the surrounding text ever contains XML-special characters.
$tailAfterPhraseXml = $tailAfterPhrase.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

$newXml = "<pkg:package xmlns:pkg=""http://schemas.microsoft.com/office/2006/xmlPackage"">" + `
    "<pkg:part pkg:name=""/word/document.xml"" pkg:contentType=""application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"">" + `
    "<pkg:xmlData><w:document xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main"">" + `
    "<w:body><w:p>" + `
    "<w:r><w:t xml:space=""preserve""> T</w:t></w:r>" + `
    "<w:r><w:t>o fix this issue</w:t></w:r>" + `
    "<w:r><w:t>,</w:t></w:r>" + `
    "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" + `
    "<w:r><w:t>I</w:t></w:r>" + `
    "<w:r><w:t xml:space=""preserve"">" + " created" + $tailAfterPhraseXml + "</w:t></w:r>" + `
    "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$editRange.InsertXML($newXml)
